# Apply 'Updated symbol list' crypto price refresh to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.20'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '16'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '25.04'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '16'
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.181'
$ws.Range("E4").Value = '3HuobiTokenHT'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '16'
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05658'
$ws.Range("E5").Value = '4CronosCRO'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '16'
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.523'
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '16'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2.965'
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '16'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8124'
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '16'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8375'
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '16'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1331'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '16'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06947'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '16'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02837'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '16'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09411'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '16'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001524'
$ws.Range("E14").Value = '13BitForexTokenBF'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '16'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005974'
$ws.Range("E15").Value = '14OneONE'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '16'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006155'
$ws.Range("E16").Value = '15TigerCashTCH'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '16'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.503'
$ws.Range("E17").Value = '16LEOLEO'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '16'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '16'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3164'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '16'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03214'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '16'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '16'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.764'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '16'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04706'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '16'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1369'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '16'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001237'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '16'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004540'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '16'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009690'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '16'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001951'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '16'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '16'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '16'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '16'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '16'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '16'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '16'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '16'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '16'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '16'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '16'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '16'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '16'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1052'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '16'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.002697'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '16'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003276'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '16'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007385'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '16'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005274'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '16'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '16'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.2198'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '16'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '16'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '16'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '16'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '16'
